# Penalty/Reward System changes: shift Week_Start_Date forward by one week
# and update MyForecast values on the "Forecast Comparison" sheet; refresh
# the derived metrics on the "Summary" sheet accordingly.

$wb = $excel.ActiveWorkbook

# ---- Forecast Comparison sheet ----
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

function Set-TextValue {
    param($range, $value)
    # Force the cell to stay a text cell (the source values are literal
    # strings, e.g. ISO dates, not real Excel dates/numbers) and make sure
    # no leftover formatting sticks to the cell once the value is in place.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$weekUpdates = @(
    @{ Row = 2;  Date = "2025-01-12"; Forecast = 62 },
    @{ Row = 3;  Date = "2025-01-19"; Forecast = 74 },
    @{ Row = 4;  Date = "2025-01-26"; Forecast = 75 },
    @{ Row = 5;  Date = "2025-02-02"; Forecast = 67 },
    @{ Row = 6;  Date = "2025-02-09"; Forecast = 62 },
    @{ Row = 7;  Date = "2025-02-16"; Forecast = 69 },
    @{ Row = 8;  Date = "2025-02-23"; Forecast = 82 },
    @{ Row = 9;  Date = "2025-03-02"; Forecast = 85 },
    @{ Row = 10; Date = "2025-03-09"; Forecast = 77 },
    @{ Row = 11; Date = "2025-03-16"; Forecast = 70 },
    @{ Row = 12; Date = "2025-03-23"; Forecast = 76 },
    @{ Row = 13; Date = "2025-03-30"; Forecast = 59 },
    @{ Row = 14; Date = "2025-04-06"; Forecast = 58 },
    @{ Row = 15; Date = "2025-04-13"; Forecast = 54 },
    @{ Row = 16; Date = "2025-04-20"; Forecast = 76 },
    @{ Row = 17; Date = "2025-04-27"; Forecast = 54 }
)

foreach ($update in $weekUpdates) {
    Set-TextValue $ws1.Cells.Item($update.Row, 2) $update.Date
    $ws1.Cells.Item($update.Row, 4).Value = $update.Forecast
}

# ---- Summary sheet ----
$ws2 = $wb.Worksheets.Item("Summary")

Set-TextValue $ws2.Range("B2")  "2024-05-05 to 2025-01-05"
Set-TextValue $ws2.Range("B4")  "68"
Set-TextValue $ws2.Range("B5")  "27"
Set-TextValue $ws2.Range("B8")  "991 units"
Set-TextValue $ws2.Range("B9")  "1099"
Set-TextValue $ws2.Range("B10") "576"
Set-TextValue $ws2.Range("B11") "278"
Set-TextValue $ws2.Range("B12") "85"
Set-TextValue $ws2.Range("B13") "2025-03-02"
Set-TextValue $ws2.Range("B14") "54"
